$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update quantity values (bug fixes in library calculations)
$ws.Range("B2").Value = 1.5
$ws.Range("B4").Value = 22.5
$ws.Range("B7").Value = 8
$ws.Range("B8").Value = 2

# Remove the now-obsolete "Área de Aço" row entirely
$ws.Rows.Item(10).Delete()
